$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume 1h (E) columns for each coin row.
# Price values are forced to Text via a temporary "@" NumberFormat (then
# ClearFormats restores the default/general style) so strings like "22.70"
# or "1.930.88" are preserved exactly instead of being auto-coerced to numbers.

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "29.132.69"
$cell.ClearFormats()
$ws.Range("E2").Value = "  +1.69%  "

$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "1.930.88"
$cell.ClearFormats()
$ws.Range("E3").Value = "  +2.55%  "

$ws.Range("E4").Value = "  +0.40%  "

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "327.23"
$cell.ClearFormats()
$ws.Range("E5").Value = "  +1.83%  "

$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "1.005"
$cell.ClearFormats()
$ws.Range("E6").Value = "  +0.39%  "

$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = "0.4617"
$cell.ClearFormats()
$ws.Range("E7").Value = "  +1.21%  "

$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "0.3838"
$cell.ClearFormats()
$ws.Range("E8").Value = "  +1.75%  "

$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "0.07787"
$cell.ClearFormats()
$ws.Range("E9").Value = "  +1.44%  "

$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "0.9812"
$cell.ClearFormats()
$ws.Range("E10").Value = "  +2.64%  "

$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "22.70"
$cell.ClearFormats()
$ws.Range("E11").Value = "  +3.91%  "

$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "1.965.64"
$cell.ClearFormats()
$ws.Range("E12").Value = "  +4.39%  "

$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "5.709"
$cell.ClearFormats()
$ws.Range("E13").Value = "  +1.50%  "

$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "6.986"
$cell.ClearFormats()
$ws.Range("E14").Value = "  +0.89%  "

$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "0.07072"
$cell.ClearFormats()
$ws.Range("E15").Value = "  +0.60%  "

$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "1.008"
$cell.ClearFormats()
$ws.Range("E16").Value = "  +0.39%  "

$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "84.54"
$cell.ClearFormats()
$ws.Range("E17").Value = "  +2.51%  "

$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "0.000009574"
$cell.ClearFormats()
$ws.Range("E18").Value = "  +1.13%  "

$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "16.75"
$cell.ClearFormats()
$ws.Range("E19").Value = "  +1.31%  "

$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "1.005"
$cell.ClearFormats()
$ws.Range("E20").Value = "  +0.32%  "

$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "29.154.40"
$cell.ClearFormats()
$ws.Range("E21").Value = "  +1.75%  "

$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "5.359"
$cell.ClearFormats()
$ws.Range("E22").Value = "  +1.44%  "

$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "10.99"
$cell.ClearFormats()
$ws.Range("E23").Value = "  +1.47%  "

$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "2.097"
$cell.ClearFormats()
$ws.Range("E24").Value = "  +0.56%  "

$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "158.26"
$cell.ClearFormats()
$ws.Range("E25").Value = "  +2.31%  "

$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "19.17"
$cell.ClearFormats()
$ws.Range("E26").Value = "  +1.40%  "

$ws.Range("E27").Value = "  +1.50%  "

$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "118.19"
$cell.ClearFormats()
$ws.Range("E28").Value = "  +1.45%  "

$ws.Range("E29").Value = "  +2.98%  "

$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "0.09356"
$cell.ClearFormats()
$ws.Range("E30").Value = "  +1.58%  "

$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "0.8666"
$cell.ClearFormats()
$ws.Range("E31").Value = "  +3.19%  "

$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "5.149"
$cell.ClearFormats()
$ws.Range("E32").Value = "  +2.26%  "

$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "1.252"
$cell.ClearFormats()
$ws.Range("E33").Value = "  +1.16%  "

$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "3.028"
$cell.ClearFormats()
$ws.Range("E34").Value = "  +3.86%  "

$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "0.05720"
$cell.ClearFormats()
$ws.Range("E35").Value = "  +1.45%  "

$ws.Range("E36").Value = "  +1.46%  "

$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "1.006"
$cell.ClearFormats()
$ws.Range("E37").Value = "  +0.34%  "

$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "0.02060"
$cell.ClearFormats()
$ws.Range("E38").Value = "  +2.32%  "

$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "3.089"
$cell.ClearFormats()
$ws.Range("E39").Value = "  +15.42%  "

$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "7.565"
$cell.ClearFormats()
$ws.Range("E40").Value = "  +2.08%  "

$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "0.5531"
$cell.ClearFormats()
$ws.Range("E41").Value = "  +1.55%  "

$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "0.1760"
$cell.ClearFormats()
$ws.Range("E42").Value = "  +1.21%  "

$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "9.378"
$cell.ClearFormats()
$ws.Range("E43").Value = "  +2.67%  "

$ws.Range("E44").Value = "  -4.15%  "

$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "2.224"
$cell.ClearFormats()
$ws.Range("E45").Value = "  +7.23%  "

$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "0.5215"
$cell.ClearFormats()
$ws.Range("E46").Value = "  +1.56%  "

$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "11.31"
$cell.ClearFormats()
$ws.Range("E47").Value = "  +1.66%  "

$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "0.06934"
$cell.ClearFormats()
$ws.Range("E48").Value = "  +2.74%  "

$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "1.779"
$cell.ClearFormats()
$ws.Range("E49").Value = "  +1.14%  "

$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "110.58"
$cell.ClearFormats()
$ws.Range("E50").Value = "  +0.05%  "

$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "1.006"
$cell.ClearFormats()
$ws.Range("E51").Value = "  +0.46%  "
